# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados..." timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 18 de Abril de 2020 a las 17:22"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 712399
$ws.Range("C4").Value = 2664
$ws.Range("D4").Value = 63768
$ws.Range("E4").Value = 611363
$ws.Range("F4").Value = 13536
$ws.Range("G4").Value = 114
$ws.Range("H4").Value = 37268

# Row 8 - Alemania
$ws.Range("B8").Value = 142325
$ws.Range("C8").Value = 928
$ws.Range("E8").Value = 52522

# Row 16 - Canada
$ws.Range("B16").Value = 32412
$ws.Range("C16").Value = 485
$ws.Range("E16").Value = 20523
$ws.Range("G16").Value = 36
$ws.Range("H16").Value = 1346

# Row 31 - Rumania
$ws.Range("E31").Value = 6267
$ws.Range("G31").Value = 10
$ws.Range("H31").Value = 421

# Row 43 - Singapur
$ws.Range("D43").Value = 740
$ws.Range("E43").Value = 5241
$ws.Range("F43").Value = 23

# Row 60 - Grecia
$ws.Range("B60").Value = 2235
$ws.Range("C60").Value = 11
$ws.Range("E60").Value = 1856
$ws.Range("F60").Value = 67
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 110

# Row 84 - Bulgaria
$ws.Range("B84").Value = 878
$ws.Range("C84").Value = 32
$ws.Range("E84").Value = 684
$ws.Range("F84").Value = 35

# Row 87 - Republica de Chipre
$ws.Range("B87").Value = 761
$ws.Range("C87").Value = 11
$ws.Range("E87").Value = 672

# Row 92 - Principado de Andorra
$ws.Range("B92").Value = 704
$ws.Range("C92").Value = 8
$ws.Range("D92").Value = 205
$ws.Range("E92").Value = 464

# Row 104 - San Marino
$ws.Range("B104").Value = 455
$ws.Range("C104").Value = 20
$ws.Range("D104").Value = 60
$ws.Range("E104").Value = 356
